$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric to Excel; force them to remain plain text
# so they match the source workbook (which stores these as literal strings).
$textCells = @("D5", "D6", "D8", "D10", "D11", "D12", "D13", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "57.864.95"
$ws.Range("E2").Value = "  -3.99%  "
$ws.Range("D3").Value = "2.953.24"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "558.57"
$ws.Range("E5").Value = "  -2.69%  "
$ws.Range("D6").Value = "131.56"
$ws.Range("E6").Value = "  +5.21%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  +3.45%  "
$ws.Range("D9").Value = "2.948.30"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").Value = "0.127"
$ws.Range("E10").Value = "  -2.95%  "
$ws.Range("D11").Value = "4.78"
$ws.Range("E11").Value = "  -6.09%  "
$ws.Range("D12").Value = "0.445"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("D13").Value = "0.0000222"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").Value = "33.04"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "3.443.16"
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("D17").Value = "6.80"
$ws.Range("E17").Value = "  +6.12%  "
$ws.Range("D18").Value = "2.951.32"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").Value = "57.841.96"
$ws.Range("E19").Value = "  -4.14%  "
$ws.Range("D20").Value = "418.39"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").Value = "13.18"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").Value = "0.685"
$ws.Range("E22").Value = "  +3.10%  "
$ws.Range("D23").Value = "6.97"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "13.12"
$ws.Range("E24").Value = "  +2.62%  "
$ws.Range("D25").Value = "79.84"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "2.49"
$ws.Range("E28").Value = "  -2.75%  "
$ws.Range("D29").Value = "7.55"
$ws.Range("E29").Value = "  +3.90%  "
$ws.Range("D30").Value = "2.01"
$ws.Range("E30").Value = "  +5.33%  "
$ws.Range("D31").Value = "25.23"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").Value = "5.98"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").Value = "0.0974"
$ws.Range("E33").Value = "  +5.27%  "
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").Value = "0.951"
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "5.64"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("B36").Value = "Stacks"
$ws.Range("C36").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D36").Value = "2.08"
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "48.46"
$ws.Range("E37").Value = "  -3.31%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0682"
$ws.Range("E38").Value = "  +4.60%  "
$ws.Range("D39").Value = "8.68"
$ws.Range("E39").Value = "  +3.24%  "
$ws.Range("D40").Value = "2.56"
$ws.Range("E40").Value = "  +4.05%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.107"
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0344"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.694.17"
$ws.Range("E43").Value = "  +1.61%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "372.42"
$ws.Range("E44").Value = "  -2.78%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "123.70"
$ws.Range("E46").Value = "  +3.21%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").Value = "0.238"
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("D48").Value = "0.110"
$ws.Range("E48").Value = "  +2.76%  "
$ws.Range("D49").Value = "1.96"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").Value = "23.23"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "2.00"
$ws.Range("E51").Value = "  +0.77%  "
